$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J25").Value = 1.974821465592385
$ws.Range("K25").Value = -1.411950589942853
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = -1.4
$ws.Range("H27").Value = 2.014950089106234
$ws.Range("I27").Value = -1.36604759998699
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -1.4
$ws.Range("F29").Value = 1.989179182172823
$ws.Range("G29").Value = -1.424065704916122
$ws.Range("H29").Value = -1.414809667530936
$ws.Range("I29").Value = -0.07315814057027888
$ws.Range("J29").Value = 0.1125927754601833
$ws.Range("K29").Value = -0.0679294223966645
$ws.Range("E30").Value = 1.998395424983189
$ws.Range("F30").Value = -1.406951595395129
$ws.Range("G30").Value = -1.449441663323289
$ws.Range("H30").Value = -0.1271799142555619
$ws.Range("I30").Value = 0.09173371707086646
$ws.Range("J30").Value = -0.1309495813748543
$ws.Range("D31").Value = 1.963093084223642
$ws.Range("E31").Value = -1.431606654211006
$ws.Range("F31").Value = -1.321860256330814
$ws.Range("G31").Value = -0.1022580996764944
$ws.Range("H31").Value = 0.09641112172120742
$ws.Range("I31").Value = -0.1271916431867545
$ws.Range("C32").Value = 1.4
$ws.Range("D32").Value = -2.4
$ws.Range("E32").Value = -1.4
$ws.Range("F32").Value = -0.1
$ws.Range("G32").Value = 0.1
$ws.Range("H32").Value = -0.1
$ws.Range("B33").Value = 1.555942877607013
$ws.Range("C33").Value = -1.436405885542493
$ws.Range("D33").Value = -1.391705137478664
$ws.Range("E33").Value = -0.1489278717582538
$ws.Range("F33").Value = 0.05372211536164659
$ws.Range("G33").Value = -0.1106874806568095
$ws.Range("H33").Value = -1.214579810576154
$ws.Range("I33").Value = -0.1090444662145961
$ws.Range("J33").Value = 0.06217141579022616
$ws.Range("K33").Value = -0.3372075916540097
$ws.Range("B34").Value = -2.125594520446883
$ws.Range("C34").Value = -1.193060812328801
$ws.Range("D34").Value = 0.2243356571731152
$ws.Range("E34").Value = 0.1386281808446615
$ws.Range("F34").Value = -0.05202793475815978
$ws.Range("G34").Value = -1.156495754450727
$ws.Range("H34").Value = -0.05060148459607872
$ws.Range("I34").Value = 0.1307542111166967
$ws.Range("J34").Value = -0.2548116829487251
$ws.Range("B35").Value = -1.929266060900714
$ws.Range("C35").Value = -0.3499020898492627
$ws.Range("D35").Value = -0.201987220327034
$ws.Range("E35").Value = -0.3435822669317387
$ws.Range("F35").Value = -1.353809024284176
$ws.Range("G35").Value = -0.1057731767660739
$ws.Range("H35").Value = 0.1735898665631708
$ws.Range("I35").Value = -0.2498049437364732
$ws.Range("B36").Value = -0.2550786956675604
$ws.Range("C36").Value = 0.3188137050645766
$ws.Range("D36").Value = -0.07489936260599347
$ws.Range("E36").Value = -1.155692413752599
$ws.Range("F36").Value = -0.06872572511066544
$ws.Range("G36").Value = 0.09229914305540476
$ws.Range("H36").Value = -0.2986516261125417
$ws.Range("B37").Value = -2.068330733759602
$ws.Range("C37").Value = -0.05755194183036486
$ws.Range("D37").Value = -0.6563389999335666
$ws.Range("E37").Value = 0.2226642396287572
$ws.Range("F37").Value = 0.2953926851715814
$ws.Range("G37").Value = -0.2642947319481946
$ws.Range("H37").Value = -0.006178249755636078
$ws.Range("I37").Value = -0.4137278119927412
$ws.Range("J37").Value = 0.3525807978017975
$ws.Range("K37").Value = 0.05185204303082339
$ws.Range("B38").Value = -0.1398269496361429
$ws.Range("C38").Value = -1.169461031008741
$ws.Range("D38").Value = -0.5943583309198688
$ws.Range("E38").Value = 0.06854235289222463
$ws.Range("F38").Value = -0.3481277478808779
$ws.Range("G38").Value = -0.02337583031178853
$ws.Range("H38").Value = -0.4413128671451602
$ws.Range("I38").Value = 0.3577964719108673
$ws.Range("J38").Value = 0.05716851371525165
$ws.Range("B39").Value = -0.9330598254158777
$ws.Range("C39").Value = -0.4460299254023569
$ws.Range("D39").Value = 0.08202667586031168
$ws.Range("E39").Value = -0.3391619959180498
$ws.Range("F39").Value = -0.04717778368150999
$ws.Range("G39").Value = -0.4455317243747756
$ws.Range("H39").Value = 0.3530597661666662
$ws.Range("I39").Value = 0.05580541522540706
$ws.Range("B40").Value = -0.06208512593411045
$ws.Range("C40").Value = 0.1491238821950708
$ws.Range("D40").Value = -0.2912725171533643
$ws.Range("E40").Value = 0.0166965182527673
$ws.Range("F40").Value = -0.4204542095353839
$ws.Range("G40").Value = 0.4428244917343195
$ws.Range("H40").Value = 0.1115427723019478
$ws.Range("B41").Value = 0.0543279172719193
$ws.Range("C41").Value = -0.2660122439511319
$ws.Range("D41").Value = 0.03050620287519542
$ws.Range("E41").Value = -0.3940545263449339
$ws.Range("F41").Value = 0.3743020454662089
$ws.Range("G41").Value = 0.1430102132627523
$ws.Range("H41").Value = -0.04414965315436631
$ws.Range("I41").Value = -0.3956711847722491
$ws.Range("J41").Value = -0.7708214880419186
$ws.Range("K41").Value = -0.7793630692420719
$ws.Range("B42").Value = -0.8154625125417774
$ws.Range("C42").Value = 0.04047683172432792
$ws.Range("D42").Value = -0.4438744219204829
$ws.Range("E42").Value = 0.4182139657863181
$ws.Range("F42").Value = 0.1270960776614676
$ws.Range("G42").Value = 0.003949865091989385
$ws.Range("H42").Value = -0.409960147238442
$ws.Range("I42").Value = -0.7920769438047226
$ws.Range("J42").Value = -0.7739584851833533
$ws.Range("B43").Value = -0.1813602613933202
$ws.Range("C43").Value = -0.3563049886777505
$ws.Range("D43").Value = 0.4297390449641996
$ws.Range("E43").Value = 0.1117284986505434
$ws.Range("F43").Value = 0.04171417100533775
$ws.Range("G43").Value = -0.3576160308555183
$ws.Range("H43").Value = -0.7589389891227414
$ws.Range("I43").Value = -0.8063269205198546
$ws.Range("B44").Value = -0.4148081973238454
$ws.Range("C44").Value = 0.4329735507955001
$ws.Range("D44").Value = 0.06331082616568086
$ws.Range("E44").Value = -0.02833738387902709
$ws.Range("F44").Value = -0.4111121799995046
$ws.Range("G44").Value = -0.8175768037466752
$ws.Range("H44").Value = -0.8328508019613353
$ws.Range("B45").Value = 0.4291840095081929
$ws.Range("C45").Value = 0.1297646319698398
$ws.Range("D45").Value = 0.0425663075747533
$ws.Range("E45").Value = -0.363896397443659
$ws.Range("F45").Value = -0.7552073568718886
$ws.Range("G45").Value = -0.8431584147299345
$ws.Range("H45").Value = 0.877766213792631
$ws.Range("I45").Value = -0.4745070531313464
$ws.Range("J45").Value = ""
$ws.Range("B46").Value = 0.06678455670716371
$ws.Range("C46").Value = 0.02006819488540634
$ws.Range("D46").Value = -0.3815327228478058
$ws.Range("E46").Value = -0.8053621644989952
$ws.Range("F46").Value = -0.8311992292910589
$ws.Range("G46").Value = 0.902157651914689
$ws.Range("H46").Value = -0.5424296176071748
$ws.Range("I46").Value = ""
$ws.Range("B47").Value = 0.00001303303454188581
$ws.Range("C47").Value = -0.3651743008299998
$ws.Range("D47").Value = -0.7715150510102097
$ws.Range("E47").Value = -0.7759485756677198
$ws.Range("F47").Value = 0.92753124260505
$ws.Range("G47").Value = -0.5192183799751632
$ws.Range("H47").Value = ""
$ws.Range("B48").Value = -0.406125572440377
$ws.Range("C48").Value = -0.7784671131096405
$ws.Range("D48").Value = -0.8492888959252406
$ws.Range("E48").Value = 0.852940314980188
$ws.Range("F48").Value = -0.4987500681587066
$ws.Range("G48").Value = ""
$ws.Range("B49").Value = -0.7512006267496926
$ws.Range("C49").Value = -0.755887583918243
$ws.Range("D49").Value = 0.8513030061134607
$ws.Range("E49").Value = -0.539931685282612
$ws.Range("F49").Value = ""
$ws.Range("B50").Value = -0.7522304086392605
$ws.Range("C50").Value = 0.946789236422905
$ws.Range("D50").Value = -0.533712758908308
$ws.Range("E50").Value = ""
$ws.Range("B51").Value = 1.371451352842971
$ws.Range("C51").Value = -0.4529753548703491
$ws.Range("D51").Value = ""
$ws.Range("B52").Value = -0.526054543893956
$ws.Range("C52").Value = ""
$ws.Range("B53").Value = ""
